$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet tab ---
$ws.Name = "salads_greenmountain"

# --- Resize the existing table to make room for the new "LeaveEmpty" column ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G4"))

# --- Cell content updates ---
# New header for the freshly added column G
$ws.Range("G1").Value = "LeaveEmpty"

# nutritionLabel links now point at per-item pages instead of the shared "waffles" placeholder
$ws.Range("F2").Value = "cobbsalad"
$ws.Range("F3").Value = "kalecaesar"
$ws.Range("F4").Value = "housesalad"

# House Salad now carries the standard allergen note (previously blank)
$ws.Range("C4").Value = "No known priority allergens"

# --- Table / column naming cleanup ---
$lo.Name = "Table32"
$lo.ListColumns.Item(7).Name = "LeaveEmpty"

# --- Minor column width tweak on column F ---
$ws.Columns.Item(6).ColumnWidth = 21.83

# --- Restore the cursor/selection position left by the editing session ---
$ws.Range("C16").Select()
